# Data update and directory clean up
# Extend the tracked-days data (columns L:P for rows 2-28, columns J:N for rows 29-36)
# with newly recorded observations. Default value is "NA"; a handful of cups were
# actually released on a given day and get the "released" status instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-28: new data lives in columns L through P
$colsLP = @("L", "M", "N", "O", "P")
for ($r = 2; $r -le 28; $r++) {
    foreach ($col in $colsLP) {
        $ws.Range("$col$r").Value = "NA"
    }
}

# Rows 29-36: new data lives in columns J through N
$colsJN = @("J", "K", "L", "M", "N")
for ($r = 29; $r -le 36; $r++) {
    foreach ($col in $colsJN) {
        $ws.Range("$col$r").Value = "NA"
    }
}

# Specific cups that were released rather than still "NA"
$ws.Range("L7").Value = "released"
$ws.Range("M11").Value = "released"
$ws.Range("M13").Value = "released"
$ws.Range("J31").Value = "released"
$ws.Range("L33").Value = "released"
$ws.Range("M33").Value = "released"
$ws.Range("J34").Value = "released"

# Restore the scroll position / selection left by the author when they saved
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("O31").Select()
